# Update "人气"/count column (F) values on the "展览" sheet and the
# aggregated "全部类型" sheet, reflecting newer figures pulled at build time.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 29
$wsExhibition.Range("F5").Value = 4733
$wsExhibition.Range("F6").Value = 167
$wsExhibition.Range("F7").Value = 59
$wsExhibition.Range("F9").Value = 35

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F8").Value = 29
$wsAll.Range("F9").Value = 4733
$wsAll.Range("F10").Value = 167
$wsAll.Range("F11").Value = 59
$wsAll.Range("F14").Value = 35
